# [FIX] fixation of biomass reaction flux at 0.11
#
# The underlying flux-variability simulation was rerun with the biomass
# reaction flux fixed at 0.11. This changed most of the reported
# flux/minimum/maximum values, and the
# "L-Glutamate__Phellogen_Day_sp_exchange" row (previously row 14, a
# non-robust outlier) is no longer part of the result set, so that whole
# row is removed and everything below it shifts up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for L-Glutamate__Phellogen_Day_sp_exchange; rows 15-18
# shift up to become rows 14-17 and the sheet dimension shrinks to A1:E17.
$ws.Rows.Item(14).Delete()

# Refresh the flux / minimum / maximum / robust values for every remaining
# data row (2-17) to match the rerun results.

$ws.Cells.Item(2, 2).Value = -0.02371014509725695
$ws.Cells.Item(2, 3).Value = -0.02800043639600445
$ws.Cells.Item(2, 4).Value = 0.02565806070686093
$ws.Cells.Item(2, 5).Value = $False

$ws.Cells.Item(3, 2).Value = -0.2304245655886423
$ws.Cells.Item(3, 3).Value = -0.2304310523523641
$ws.Cells.Item(3, 4).Value = -0.2304103386042093
$ws.Cells.Item(3, 5).Value = $True

$ws.Cells.Item(4, 2).Value = -0.00300718951407934
$ws.Cells.Item(4, 3).Value = -0.003212557730932339
$ws.Cells.Item(4, 4).Value = -0.002953699178742782
$ws.Cells.Item(4, 5).Value = $True

$ws.Cells.Item(5, 2).Value = -0.0057703036764802
$ws.Cells.Item(5, 3).Value = -0.006998398981867556
$ws.Cells.Item(5, 4).Value = -0.004534755817959456
$ws.Cells.Item(5, 5).Value = $True

$ws.Cells.Item(6, 2).Value = -0.007653100707315835
$ws.Cells.Item(6, 3).Value = -0.007658336329263271
$ws.Cells.Item(6, 4).Value = -0.007632092013061378
$ws.Cells.Item(6, 5).Value = $True

$ws.Cells.Item(7, 2).Value = -0.009179095652298
$ws.Cells.Item(7, 3).Value = -0.009180723707159286
$ws.Cells.Item(7, 4).Value = -0.009168271601245401
$ws.Cells.Item(7, 5).Value = $True

$ws.Cells.Item(8, 2).Value = -0.0038465916288536
$ws.Cells.Item(8, 3).Value = -0.004665261220309521
$ws.Cells.Item(8, 4).Value = -0.003022710130898713
$ws.Cells.Item(8, 5).Value = $True

$ws.Cells.Item(9, 2).Value = -0.005717341279690275
$ws.Cells.Item(9, 3).Value = -0.00571734128654498
$ws.Cells.Item(9, 4).Value = -0.005717341263579607
$ws.Cells.Item(9, 5).Value = $True

$ws.Cells.Item(10, 2).Value = -0.006613269852417414
$ws.Cells.Item(10, 3).Value = -0.01125890894495702
$ws.Cells.Item(10, 4).Value = -0.001967630896746799
$ws.Cells.Item(10, 5).Value = $True

$ws.Cells.Item(11, 2).Value = -0.1980226309911896
$ws.Cells.Item(11, 3).Value = -0.1980420776586849
$ws.Cells.Item(11, 4).Value = -0.01521694661392164
$ws.Cells.Item(11, 5).Value = $True

$ws.Cells.Item(12, 2).Value = -0.01806263238071845
$ws.Cells.Item(12, 3).Value = -0.01809652915349802
$ws.Cells.Item(12, 4).Value = -0.01805601071758732
$ws.Cells.Item(12, 5).Value = $True

$ws.Cells.Item(13, 2).Value = -0.3459488551093048
$ws.Cells.Item(13, 3).Value = -0.3461064285972724
$ws.Cells.Item(13, 4).Value = -0.3459281046215601
$ws.Cells.Item(13, 5).Value = $True

# (was row 15) L-Asparagine__Leaf_Day_sp_exchange
$ws.Cells.Item(14, 2).Value = -0.0060052727868705
$ws.Cells.Item(14, 3).Value = -0.006017259009205319
$ws.Cells.Item(14, 4).Value = -0.005993409558369541
$ws.Cells.Item(14, 5).Value = $True

# (was row 16) L-Asparagine__Ibark_Day_sp_exchange
$ws.Cells.Item(15, 2).Value = -0.0025165694863606
$ws.Cells.Item(15, 3).Value = -0.003052170600647712
$ws.Cells.Item(15, 4).Value = -0.001980968368070639
$ws.Cells.Item(15, 5).Value = $True

# (was row 17) L-Serine__Leaf_Day_sp_exchange
$ws.Cells.Item(16, 2).Value = -0.1002073721757014
$ws.Cells.Item(16, 3).Value = -0.2755390577924466
$ws.Cells.Item(16, 4).Value = -0.04719036363783607
$ws.Cells.Item(16, 5).Value = $True

# (was row 18) Citrate__Leaf_Day_sp_exchange
$ws.Cells.Item(17, 2).Value = -0.007301385095231733
$ws.Cells.Item(17, 3).Value = -0.007317797398071424
$ws.Cells.Item(17, 4).Value = -0.005635533610203441
$ws.Cells.Item(17, 5).Value = $True
